{"js": "// edit.js - Office.js (Word JavaScript API) script\n// Applies two content changes to the document:\n//   1. Appends \"  (This is a change \u2013 Version for branch alternate)\" to the\n//      end of the first paragraph (\"This is a Microsoft word document.\"),\n//      keeping the leading two spaces in the default (black) color and\n//      coloring the parenthetical text dark red (#C00000).\n//   2. Adds a new, completely empty paragraph right before the end of the\n//      document (after the last existing paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Extend the first paragraph with the \"branch alternate\" annotation.\n// ---------------------------------------------------------------------\nconst firstParagraph = paragraphs.items[0];\n\n// Plain (uncolored) two spaces, appended to the existing run.\nfirstParagraph.insertText(\"  \", Word.InsertLocation.end);\n\n// Colored parenthetical text, inserted as its own run afterwards.\nconst coloredRange = firstParagraph.insertText(\n  \"(This is a change \\u2013 Version for branch alternate)\",\n  Word.InsertLocation.end\n);\ncoloredRange.font.color = \"#C00000\";\n\n// ---------------------------------------------------------------------\n// 2) Append a brand-new, empty paragraph at the very end of the document.\n// ---------------------------------------------------------------------\nconst newParagraph = body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n\n// The freshly inserted paragraph picked up the previous paragraph's style\n// (and an empty placeholder run). Replace its contents with a pristine,\n// attribute-free <w:p/> so it matches a truly blank paragraph.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p/></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst newParagraphRange = newParagraph.getRange();\nnewParagraphRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# edit.ps1 - Word COM interop script\n# Applies two content changes to the document:\n#   1. Appends \"  (This is a change \u2013 Version for branch alternate)\" to the\n#      end of the first paragraph (\"This is a Microsoft word document.\"),\n#      keeping the leading two spaces in the default (black) color and\n#      coloring the parenthetical text dark red (C00000).\n#   2. Adds a new, completely empty paragraph right before the end of the\n#      document (after the last existing paragraph).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Extend the first paragraph with the \"branch alternate\" annotation.\n# ---------------------------------------------------------------------\n$p1 = $d.Paragraphs(1)\n$beforeEnd = $p1.Range.End\n\n$newText = \"  (This is a change \" + [char]0x2013 + \" Version for branch alternate)\"\n$p1.Range.InsertAfter($newText)\n$afterEnd = $p1.Range.End\n\n# Paragraph.Range.End includes the trailing paragraph mark, so subtract 1\n# from both endpoints to get plain character offsets into the story, then\n# skip the leading two (uncolored) spaces before coloring the rest.\n$colorStart = ($beforeEnd - 1) + 2\n$colorEnd = $afterEnd - 1\n\n$colorRange = $d.Range($colorStart, $colorEnd)\n# Dark red C00000 == decimal 192 (Word colors are stored as 0x00BBGGRR, and\n# G == B == 0 here, so the BGR-encoded value equals the plain red channel).\n$colorRange.Font.Color = 192\n\n# ---------------------------------------------------------------------\n# 2) Append a brand-new, empty paragraph at the very end of the document.\n# ---------------------------------------------------------------------\n$lastIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($lastIndex)\n$lastPara.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph picked up the previous paragraph's style\n# (and an empty placeholder run). Replace its contents with a pristine,\n# attribute-free <w:p/> so it matches a truly blank paragraph.\n$newIndex = $d.Paragraphs.Count\n$newPara = $d.Paragraphs($newIndex)\n$newPara.Range.InsertXML(\"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>\")\n\nWrite-Output \"done\"\n"}
